$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F14").ClearContents()

$ws.Range("E17").Value = 7
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 1

[void]$ws.Range("F17").Select()
